$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as text in the source data; format as text
# first so Excel does not silently convert the numeric-looking strings to numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '248.86'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '22.69'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.320'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05684'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.365'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9150'
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0005852'
$ws.Range('E10').Value = '9OneONE'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1406'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07445'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03101'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.03026'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09379'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.893'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.001579'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.04806'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('B19').Value = 'UpBots'
$ws.Range('C19').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.01828'
$ws.Range('E19').Value = '18UpBotsUBXTBestin24h'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.006446'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.004992'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0009996'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0001501'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.694'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.196'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1312'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04000'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.003039'
$ws.Range('E41').Value = '40KickTokenKICKWorstin24h'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1071'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002741'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007961'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005677'
$ws.Range('E47').Value = '46CoinbaseStockTokenCOIN'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.2106'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002101'
